# Refactor the "plasmids" mock-db worksheet to support plasmid features
# (resistance genes + origins of replication) instead of the old
# single "Name" column.
#
# Before:
#   A1 = "Name"
#   A2 = "p2"
#
# After:
#   A1 = "Resistance"   B1 = "Origin"
#   (row 2 left blank)
#   A3 = "AmpR,TetR"    B3 = "p15A"
#
# The old column A ("Name"/"p2") is removed entirely (an entire-column
# delete), which is why every remaining/used column shifts one slot to
# the left, and then the new header + sample data is written in, with
# the sample row landing on row 3 (row 2 is left empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Name" column outright - this shifts column widths
# etc. left by one, matching the target layout.
$ws.Columns.Item(1).Delete()

# New header row.
$ws.Range("A1").Value = "Resistance"
$ws.Range("B1").Value = "Origin"

# New sample data, on row 3 (row 2 stays empty).
$ws.Range("A3").Value = "AmpR,TetR"
$ws.Range("B3").Value = "p15A"

# Match the refreshed selection (top-left cell instead of the old C6).
$ws.Range("A1").Select()
